$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fix: Q2 (Chargesheets submitted ... during the Year, for row "Murder with Rape/Gang Rape")
# was mistakenly showing 208; correct figure is 1.
$ws.Range("Q2").Value = 1

# --- Data fix: the "Index" column (A) was re-using the same 16-value cycle
# (1,2,3,4,5,6,7,8,15,18,21,24,27,30,34,45) for every one of the 3 years of
# crime-head data instead of running as a continuous index down the sheet.
# Renumber it sequentially (row 2 -> 1, row 3 -> 2, ... row 49 -> 48).
for ($r = 2; $r -le 49; $r++) {
    $ws.Range("A$r").Value = $r - 1
}

# --- View state: zoom to 85% and move the selection/active cell to K2
# (previously the window was scrolled to show column W onward with Z48 selected).
$excel.ActiveWindow.Zoom = 85
$ws.Range("K2").Select()
